# Fix Training Data Issue (#48)
# Data was taken from 1 day off due to way NBA stats were shown.
# Corrects the BF "Date" column (was stored as "1-19-2020-21", should be
# "2021-01-19") and a handful of stat values that were off by one day.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Date column (BF) for rows 2-31: "1-19-2020-21" -> "2021-01-19" -----
# Briefly mark the cell as Text so Excel doesn't reinterpret the
# ISO-formatted string as a date serial number, then restore the cell's
# normal (default) style so no residual number formatting is left behind.
for ($r = 2; $r -le 31; $r++) {
    $cell = $ws.Cells.Item($r, 58)
    $cell.NumberFormat = "@"
    $cell.Value = "2021-01-19"
    $cell.Style = "Normal"
}

# --- Individual stat corrections ----------------------------------------
$ws.Range("Y4").Value = 3.9
$ws.Range("AN4").Value = 2

$ws.Range("AX7").Value = 17

$ws.Range("AY8").Value = 6

$ws.Range("S10").Value = 31.7
$ws.Range("T10").Value = 42.8
$ws.Range("V10").Value = 14
$ws.Range("AV10").Value = 9

$ws.Range("AX14").Value = 19

$ws.Range("J16").Value = 92.2
$ws.Range("R16").Value = 10.6
$ws.Range("T16").Value = 45.4
$ws.Range("AR16").Value = 12
$ws.Range("AX16").Value = 17

$ws.Range("M18").Value = 39.6
$ws.Range("X18").Value = 4.8
$ws.Range("AN18").Value = 3
$ws.Range("AX18").Value = 19

$ws.Range("AR21").Value = 13

$ws.Range("AV27").Value = 8
